# Update "想去人数" (F column) figures on the 展览, 演出 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    5  = 2909
    9  = 1541
    11 = 832
    12 = 95
    13 = 2604
    15 = 1452
    16 = 6791
    18 = 6245
    19 = 7
    20 = 2198
    21 = 3018
    22 = 3418
    24 = 46
    25 = 1715
    26 = 68
    27 = 283
    28 = 858
    30 = 19
    31 = 355
    32 = 1069
    33 = 2327
    34 = 9
    36 = 332
    37 = 894
    38 = 187
    39 = 424
    40 = 484
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    20 = 8
    22 = 58
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    10 = 2909
    13 = 1541
    15 = 832
    16 = 95
    18 = 2604
    19 = 1452
    24 = 6791
    26 = 6245
    27 = 2198
    28 = 3018
    29 = 3418
    33 = 1715
    34 = 8
    36 = 283
    37 = 858
    38 = 19
    39 = 355
    40 = 58
    41 = 2327
    42 = 9
    45 = 332
    46 = 894
    47 = 187
    48 = 424
    49 = 484
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

$wb.Save()
